# Update node-name identifiers in the "Training Results" sheet.
# These cells contain free-text error messages whose embedded onnx node
# ids changed between test runs; only the numeric ids need updating.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

$ws.Range("C16").Value = "C:\Users\COCO\onnxruntime_training_cuda_python\orttraining\orttraining\python\orttraining_pybind_state.cc:621 onnxruntime::python::addObjectMethodsForTraining::<lambda_6dd399ad6691adab5d0e0423ed8ce22d>::operator () [ONNXRuntimeError] : 1 : FAIL : Type Error: Type parameter (T) of Optype (Sub) bound to different types (tensor(float) and tensor(double) in node (onnx::Pow::24062_Grad/Sub_1).`n"

$ws.Range("C39").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24181): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C51").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24271): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C52").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24273): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C64").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24326): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C65").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24328): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C69").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24346): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C70").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24348): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C93").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::24444): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C102").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::24481): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C223").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::24931): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C239").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::24936): X typestr: T, has unsupported type: tensor(uint8)"
